$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.Formula = '="' + $val + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

Set-TextValue "D2" "28.594.61"
$ws.Range("E2").Value = "  +1.60%  "
Set-TextValue "D3" "1.880.74"
$ws.Range("E3").Value = "  +1.16%  "
Set-TextValue "D4" "1.012"
$ws.Range("E4").Value = "  +0.29%  "
Set-TextValue "D5" "316.60"
$ws.Range("E5").Value = "  +1.42%  "
Set-TextValue "D6" "1.010"
$ws.Range("E6").Value = "  +0.74%  "
Set-TextValue "D7" "0.5101"
$ws.Range("E7").Value = "  +0.42%  "
Set-TextValue "D8" "0.3907"
$ws.Range("E8").Value = "  +0.36%  "
Set-TextValue "D9" "0.08403"
$ws.Range("E9").Value = "  +2.13%  "
Set-TextValue "D10" "1.105"
$ws.Range("E10").Value = "  -0.32%  "
Set-TextValue "D11" "6.235"
$ws.Range("E11").Value = "  +0.54%  "
Set-TextValue "D12" "1.875.92"
$ws.Range("E12").Value = "  +1.74%  "
Set-TextValue "D13" "20.44"
$ws.Range("E13").Value = "  +1.37%  "
Set-TextValue "D14" "7.253"
$ws.Range("E14").Value = "  +1.03%  "
Set-TextValue "D15" "1.012"
$ws.Range("E15").Value = "  +0.15%  "
Set-TextValue "D16" "0.00001106"
$ws.Range("E16").Value = "  +1.18%  "
Set-TextValue "D17" "91.40"
Set-TextValue "D18" "0.06735"
$ws.Range("E18").Value = "  +0.95%  "
Set-TextValue "D19" "17.75"
$ws.Range("E19").Value = "  +1.16%  "
Set-TextValue "D20" "1.010"
$ws.Range("E20").Value = "  +0.49%  "
Set-TextValue "D21" "5.938"
$ws.Range("E21").Value = "  +0.45%  "
Set-TextValue "D22" "28.625.03"
$ws.Range("E22").Value = "  +1.73%  "
Set-TextValue "D23" "11.11"
$ws.Range("E23").Value = "  +0.77%  "
Set-TextValue "D24" "2.242"
$ws.Range("E24").Value = "  +1.06%  "
Set-TextValue "D25" "2.087.53"
$ws.Range("E25").Value = "  +1.63%  "
Set-TextValue "D26" "161.91"
$ws.Range("E26").Value = "  +1.67%  "
Set-TextValue "D27" "20.78"
$ws.Range("E27").Value = "  +1.03%  "
Set-TextValue "D28" "2.356"
$ws.Range("E28").Value = "  -1.98%  "
Set-TextValue "D29" "126.90"
$ws.Range("E29").Value = "  +0.59%  "
Set-TextValue "D30" "0.1048"
$ws.Range("E30").Value = "  -0.58%  "
Set-TextValue "D31" "1.043"
$ws.Range("E31").Value = "  +1.25%  "
Set-TextValue "D32" "5.797"
$ws.Range("E32").Value = "  -0.52%  "
Set-TextValue "D33" "3.619"
$ws.Range("E33").Value = "  +0.32%  "
Set-TextValue "D34" "0.02466"
$ws.Range("E34").Value = "  +1.92%  "
Set-TextValue "D35" "0.06561"
$ws.Range("E35").Value = "  +1.86%  "
Set-TextValue "D36" "0.2166"
$ws.Range("E36").Value = "  +0.30%  "
Set-TextValue "D37" "8.866"
$ws.Range("E37").Value = "  -1.82%  "
Set-TextValue "D38" "5.081"
$ws.Range("E38").Value = "  +3.12%  "
Set-TextValue "D39" "1.197"
$ws.Range("E39").Value = "  +1.66%  "
Set-TextValue "D40" "1.256"
$ws.Range("E40").Value = "  +0.51%  "
Set-TextValue "D41" "0.6437"
Set-TextValue "D42" "11.13"
$ws.Range("E42").Value = "  +0.72%  "
Set-TextValue "D43" "1.010"
$ws.Range("E43").Value = "  +0.72%  "
Set-TextValue "D44" "0.6049"
Set-TextValue "D45" "13.00"
$ws.Range("E45").Value = "  +0.14%  "
Set-TextValue "D46" "3.699"
$ws.Range("E46").Value = "  +1.09%  "
Set-TextValue "D47" "2.011"
$ws.Range("E47").Value = "  +0.81%  "
Set-TextValue "D48" "1.220"
$ws.Range("E48").Value = "  +1.46%  "
Set-TextValue "D49" "122.19"
$ws.Range("E49").Value = "  +1.38%  "
Set-TextValue "D50" "1.146"
$ws.Range("E50").Value = "  -9.73%  "
Set-TextValue "D51" "0.06840"
$ws.Range("E51").Value = "  -0.38%  "

$excel.CutCopyMode = 0
